# Updated cryptos list: refresh Price (D) and Volume(1h) (E) columns
# to match the latest scrape. Values that look like plain numbers are
# written with a leading apostrophe so Excel keeps them as text (matching
# the original inline-string cell type / formatting, e.g. "7.60" must not
# collapse to 7.6).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.932.89'
$ws.Range('E2').Value = '  +0.97%  '
$ws.Range('D3').Value = '3.379.52'
$ws.Range('E3').Value = '  +0.14%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '''569.31'
$ws.Range('E5').Value = '  +0.47%  '
$ws.Range('D6').Value = '''139.99'
$ws.Range('E6').Value = '  -0.38%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').Value = '''0.472'
$ws.Range('E8').Value = '  +0.10%  '
$ws.Range('D9').Value = '''7.60'
$ws.Range('E9').Value = '  +1.52%  '
$ws.Range('E10').Value = '  -1.37%  '
$ws.Range('E11').Value = '  -0.22%  '
$ws.Range('D12').Value = '3.957.37'
$ws.Range('E12').Value = '  +0.17%  '
$ws.Range('E13').Value = '  +2.05%  '
$ws.Range('D14').Value = '''27.75'
$ws.Range('E14').Value = '  -0.98%  '
$ws.Range('D15').Value = '3.380.69'
$ws.Range('E15').Value = '  +0.11%  '
$ws.Range('E16').Value = '  -0.03%  '
$ws.Range('D17').Value = '61.064.88'
$ws.Range('E17').Value = '  +0.98%  '
$ws.Range('E18').Value = '  -1.72%  '
$ws.Range('E19').Value = '  -2.02%  '
$ws.Range('D20').Value = '''8.87'
$ws.Range('E20').Value = '  -1.53%  '
$ws.Range('D21').Value = '''380.57'
$ws.Range('E21').Value = '  -1.37%  '
$ws.Range('D22').Value = '''75.40'
$ws.Range('E22').Value = '  +3.29%  '
$ws.Range('E23').Value = '  -1.15%  '
$ws.Range('E24').Value = '  -0.14%  '
$ws.Range('D25').Value = '''0.0000113'
$ws.Range('E25').Value = '  -1.72%  '
$ws.Range('D26').Value = '3.518.19'
$ws.Range('D27').Value = '''0.192'
$ws.Range('E27').Value = '  +7.26%  '
$ws.Range('E28').Value = '  +0.04%  '
$ws.Range('D29').Value = '''7.19'
$ws.Range('E29').Value = '  -2.21%  '
$ws.Range('D30').Value = '''7.93'
$ws.Range('E30').Value = '  +0.02%  '
$ws.Range('E31').Value = '  -0.36%  '
$ws.Range('E32').Value = '  -0.04%  '
$ws.Range('E33').Value = '  -3.84%  '
$ws.Range('D34').Value = '''23.16'
$ws.Range('E34').Value = '  -1.52%  '
$ws.Range('D35').Value = '''6.90'
$ws.Range('E35').Value = '  +0.06%  '
$ws.Range('D36').Value = '''166.40'
$ws.Range('E36').Value = '  -1.15%  '
$ws.Range('D37').Value = '3.414.97'
$ws.Range('E37').Value = '  +0.29%  '
$ws.Range('D38').Value = '''4.93'
$ws.Range('E38').Value = '  +0.16%  '
$ws.Range('E39').Value = '  -3.09%  '
$ws.Range('E40').Value = '  -1.05%  '
$ws.Range('D41').Value = '''25.96'
$ws.Range('E41').Value = '  -4.05%  '
$ws.Range('E43').Value = '  +0.22%  '
$ws.Range('E44').Value = '  -1.96%  '
$ws.Range('D45').Value = '''1.63'
$ws.Range('E45').Value = '  -3.20%  '
$ws.Range('E46').Value = '  -0.38%  '
$ws.Range('D47').Value = '2.425.12'
$ws.Range('E47').Value = '  -3.70%  '
$ws.Range('D48').Value = '''22.72'
$ws.Range('E48').Value = '  -1.87%  '
$ws.Range('E49').Value = '  -2.01%  '
$ws.Range('E50').Value = '  -2.96%  '
$ws.Range('E51').Value = '  +6.27%  '
